$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Simplify the sgm dictionary: rename "dexcomg6" entries to "dexcom"
$ws.Range("C7").Value = "dexcom"
$ws.Range("C8").Value = "dexcom"
$ws.Range("C9").Value = "dexcom"

# Update active cell selection to G8
$ws.Range("G8").Select()
